# Applies the feature-selection workbook correction described by the diff:
#  - Both sheets' "Feature" column (A) is re-ordered/re-labelled to the corrected order.
#  - The "final_fail" sheet has a number of per-model boolean flags (and the
#    resulting Total column) corrected for several rows.

$wb = $excel.ActiveWorkbook
$wsFail = $wb.Worksheets.Item("final_fail")
$wsGifted = $wb.Worksheets.Item("final_gifted")

# ---------------------------------------------------------------------------
# 1. Corrected Feature-name order for "final_fail" (column A, rows 2-36).
# ---------------------------------------------------------------------------
$failFeatures = @(
    "On/off campus click ratio",
    "Submissions (% of course total)",
    "Days with no interaction",
    "Clicks (% of course total)",
    "Quizzes started",
    "Assignments viewed",
    "Largest period of inactivity (h)",
    "Total time online (min)",
    "Average session duration (min)",
    "Start of Session 1 (%)",
    "Start of Session 2 (%)",
    "Resources viewed",
    "Number of days",
    "Clicks per session",
    "Clicks per day",
    "Number of clicks",
    "Start of Session 3 (%)",
    "Clicks on campus",
    "Days with no interaction (%)",
    "Clicks on folder",
    "Clicks on course",
    "Start of Session 4 (%)",
    "Start of Session 7 (%)",
    "Forum posts",
    "Files downloaded",
    "Links viewed",
    "Discussions viewed",
    "Assignments submitted",
    "Start of Session 6 (%)",
    "Number of sessions",
    "Start of Session 10 (%)",
    "Start of Session 9 (%)",
    "Start of Session 8 (%)",
    "Start of Session 5 (%)",
    "Clicks on forum"
)

for ($i = 0; $i -lt $failFeatures.Length; $i++) {
    $row = $i + 2
    $wsFail.Range("A$row").Value = $failFeatures[$i]
}

# ---------------------------------------------------------------------------
# 2. Corrected Feature-name order for "final_gifted" (column A, rows 2-36).
# ---------------------------------------------------------------------------
$giftedFeatures = @(
    "Resources viewed",
    "Clicks per session",
    "Total time online (min)",
    "Average session duration (min)",
    "On/off campus click ratio",
    "Days with no interaction",
    "Clicks (% of course total)",
    "Largest period of inactivity (h)",
    "Start of Session 4 (%)",
    "Clicks on course",
    "Start of Session 1 (%)",
    "Start of Session 3 (%)",
    "Number of days",
    "Clicks per day",
    "Number of clicks",
    "Start of Session 7 (%)",
    "Assignments viewed",
    "Days with no interaction (%)",
    "Start of Session 2 (%)",
    "Clicks on campus",
    "Submissions (% of course total)",
    "Quizzes started",
    "Links viewed",
    "Assignments submitted",
    "Clicks on folder",
    "Start of Session 10 (%)",
    "Start of Session 6 (%)",
    "Start of Session 5 (%)",
    "Discussions viewed",
    "Forum posts",
    "Files downloaded",
    "Number of sessions",
    "Start of Session 9 (%)",
    "Start of Session 8 (%)",
    "Clicks on forum"
)

for ($i = 0; $i -lt $giftedFeatures.Length; $i++) {
    $row = $i + 2
    $wsGifted.Range("A$row").Value = $giftedFeatures[$i]
}

# ---------------------------------------------------------------------------
# 3. Corrected per-model boolean flags (and derived Total) on "final_fail".
# ---------------------------------------------------------------------------
$wsFail.Range("F6").Value = $false
$wsFail.Range("J6").Value = 5

$wsFail.Range("E7").Value = $false
$wsFail.Range("J7").Value = 5

$wsFail.Range("D10").Value = $false
$wsFail.Range("F10").Value = $true

$wsFail.Range("F16").Value = $false
$wsFail.Range("J16").Value = 4

$wsFail.Range("C17").Value = $false
$wsFail.Range("J17").Value = 4

$wsFail.Range("B18").Value = $false
$wsFail.Range("E18").Value = $true

$wsFail.Range("B19").Value = $false
$wsFail.Range("C19").Value = $false
$wsFail.Range("D19").Value = $true
$wsFail.Range("E19").Value = $true

$wsFail.Range("C20").Value = $false
$wsFail.Range("J20").Value = 3

$wsFail.Range("D21").Value = $true
$wsFail.Range("E21").Value = $false
$wsFail.Range("F21").Value = $false
$wsFail.Range("J21").Value = 3

$wsFail.Range("B23").Value = $true
$wsFail.Range("C23").Value = $false
$wsFail.Range("E23").Value = $false
$wsFail.Range("F23").Value = $true

$wsFail.Range("B24").Value = $true
$wsFail.Range("C24").Value = $false
$wsFail.Range("D24").Value = $false
$wsFail.Range("F24").Value = $true

$wsFail.Range("C25").Value = $false
$wsFail.Range("J25").Value = 2

$wsFail.Range("C26").Value = $false
$wsFail.Range("D26").Value = $true
$wsFail.Range("F26").Value = $false
$wsFail.Range("J26").Value = 2

$wsFail.Range("C28").Value = $false
$wsFail.Range("D28").Value = $true

$wsFail.Range("D30").Value = $false
$wsFail.Range("F30").Value = $true

$wsFail.Range("C31").Value = $false
$wsFail.Range("J31").Value = 1

$wsFail.Range("C32").Value = $false
$wsFail.Range("J32").Value = 1

$wsFail.Range("C33").Value = $false
$wsFail.Range("J33").Value = 1

$wsFail.Range("C34").Value = $false
$wsFail.Range("J34").Value = 1

$wsFail.Range("D35").Value = $false
$wsFail.Range("J35").Value = 1
